# Add a new row (77) of data to the "inflationbreakdown" sheet, following
# the existing monthly time series (date in column A, weights in B:F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 77
$prevRow = $row - 1

# New data values (date serial 45383 = 2024-04-01)
$ws.Cells.Item($row, 1).Value = 45383
$ws.Cells.Item($row, 2).Value = 0.22839
$ws.Cells.Item($row, 3).Value = 0.1265
$ws.Cells.Item($row, 4).Value = 0.20794
$ws.Cells.Item($row, 5).Value = 0.43721
$ws.Cells.Item($row, 6).Value = 0.34975

# Match the date formatting/style used by the rest of column A (style index
# carries the custom date number format) by copying formats from the cell
# directly above, the same way a user dragging/filling the series would.
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
